$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaderNames = @("11000","12000","13000","14000","15000")
$newHeaderCols  = @("M","N","O","P","Q")

# --- Step 1: resize the 3 ListObjects (tables) so they span through column Q. ---
# --- We also "prime" each tables new ListColumns with the header names so the
# --- table part XML (tableColumn name=...) gets the real names instead of the
# --- Excel-generated defaults ("Column12".."Column16"). ---

$lo1 = $ws.ListObjects.Item(1)
$lo1.Resize($ws.Range("B3:Q7"))
for ($i = 0; $i -lt 5; $i++) {
    $col = $lo1.ListColumns.Item(12 + $i)
    $col.Range.Cells.Item(1,1).Value = $newHeaderNames[$i]
}

$lo2 = $ws.ListObjects.Item(2)
$lo2.Resize($ws.Range("B11:Q15"))
for ($i = 0; $i -lt 5; $i++) {
    $col = $lo2.ListColumns.Item(12 + $i)
    $col.Range.Cells.Item(1,1).Value = $newHeaderNames[$i]
}

$lo3 = $ws.ListObjects.Item(3)
$lo3.Resize($ws.Range("B19:Q23"))
for ($i = 0; $i -lt 5; $i++) {
    $col = $lo3.ListColumns.Item(12 + $i)
    $col.Range.Cells.Item(1,1).Value = $newHeaderNames[$i]
}

# --- Step 2: write the real header-row cells (row 3 / 11 / 19) as TEXT values,
# --- matching the original workbook where the numeric-looking headers ("1000",
# --- "2000", ...) are stored as shared strings, not numbers. We briefly force a
# --- text number format so the engine does not auto-convert "11000" -> 11000,
# --- then restore the default style so no stray "s=" attribute is left behind. ---
foreach ($hdrRow in @(3, 11, 19)) {
    for ($i = 0; $i -lt 5; $i++) {
        $addr = "$($newHeaderCols[$i])$hdrRow"
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $newHeaderNames[$i]
        $ws.Range($addr).Style = "Normal"
    }
}

# --- Step 3: fix the "subtration_*" -> "subtraction_*" typo in the B-column row labels. ---
$ws.Range("B4").Value = "addition_csr"
$ws.Range("B5").Value = "subtraction_csr"
$ws.Range("B6").Value = "addition_csc"
$ws.Range("B7").Value = "subtraction_csc"
$ws.Range("B12").Value = "addition_csr"
$ws.Range("B13").Value = "subtraction_csr"
$ws.Range("B14").Value = "addition_csc"
$ws.Range("B15").Value = "subtraction_csc"
$ws.Range("B20").Value = "addition_csr"
$ws.Range("B21").Value = "subtraction_csr"
$ws.Range("B22").Value = "addition_csc"
$ws.Range("B23").Value = "subtraction_csc"

# --- Step 4: refresh all numeric data cells C..Q for every data row (existing columns
# --- C..L get updated averages, new columns M..Q get the 11000..15000 series). ---

$ws.Range("C4").Value = 0.007923999999999999
$ws.Range("D4").Value = 0.027088
$ws.Range("E4").Value = 0.06057100000000001
$ws.Range("F4").Value = 0.113555
$ws.Range("G4").Value = 0.181902
$ws.Range("H4").Value = 0.274877
$ws.Range("I4").Value = 0.3962290000000001
$ws.Range("J4").Value = 0.539411
$ws.Range("K4").Value = 0.7195550000000001
$ws.Range("L4").Value = 0.9307650000000001
$ws.Range("M4").Value = 1.198896
$ws.Range("N4").Value = 1.445078
$ws.Range("O4").Value = 1.769522
$ws.Range("P4").Value = 2.139689
$ws.Range("Q4").Value = 2.554645

$ws.Range("C5").Value = 0.008118
$ws.Range("D5").Value = 0.027634
$ws.Range("E5").Value = 0.06273199999999998
$ws.Range("F5").Value = 0.119001
$ws.Range("G5").Value = 0.191547
$ws.Range("H5").Value = 0.287118
$ws.Range("I5").Value = 0.407289
$ws.Range("J5").Value = 0.551316
$ws.Range("K5").Value = 0.737119
$ws.Range("L5").Value = 0.953928
$ws.Range("M5").Value = 1.197954
$ws.Range("N5").Value = 1.460427
$ws.Range("O5").Value = 1.778972
$ws.Range("P5").Value = 2.166445
$ws.Range("Q5").Value = 2.594647

$ws.Range("C6").Value = 0.008089
$ws.Range("D6").Value = 0.027413
$ws.Range("E6").Value = 0.063056
$ws.Range("F6").Value = 0.118821
$ws.Range("G6").Value = 0.191288
$ws.Range("H6").Value = 0.284725
$ws.Range("I6").Value = 0.409952
$ws.Range("J6").Value = 0.556116
$ws.Range("K6").Value = 0.7298980000000002
$ws.Range("L6").Value = 0.9400960000000002
$ws.Range("M6").Value = 1.180584
$ws.Range("N6").Value = 1.456368
$ws.Range("O6").Value = 1.779277
$ws.Range("P6").Value = 2.152697
$ws.Range("Q6").Value = 2.599839

$ws.Range("C7").Value = 0.008193
$ws.Range("D7").Value = 0.028129
$ws.Range("E7").Value = 0.063433
$ws.Range("F7").Value = 0.118838
$ws.Range("G7").Value = 0.195204
$ws.Range("H7").Value = 0.291445
$ws.Range("I7").Value = 0.4118709999999999
$ws.Range("J7").Value = 0.560257
$ws.Range("K7").Value = 0.746381
$ws.Range("L7").Value = 0.953993
$ws.Range("M7").Value = 1.20354
$ws.Range("N7").Value = 1.474232
$ws.Range("O7").Value = 1.822403
$ws.Range("P7").Value = 2.207352
$ws.Range("Q7").Value = 2.658625

$ws.Range("C12").Value = 0.014075
$ws.Range("D12").Value = 0.05599
$ws.Range("E12").Value = 0.138048
$ws.Range("F12").Value = 0.271246
$ws.Range("G12").Value = 0.464772
$ws.Range("H12").Value = 0.730065
$ws.Range("I12").Value = 1.081819
$ws.Range("J12").Value = 1.522875
$ws.Range("K12").Value = 2.07307
$ws.Range("L12").Value = 2.735982
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = 4.342772
$ws.Range("O12").Value = 5.358774
$ws.Range("P12").Value = 6.547840000000001
$ws.Range("Q12").Value = 7.883540999999999

$ws.Range("C13").Value = 0.014366
$ws.Range("D13").Value = 0.058123
$ws.Range("E13").Value = 0.142464
$ws.Range("F13").Value = 0.276263
$ws.Range("G13").Value = 0.476504
$ws.Range("H13").Value = 0.738939
$ws.Range("I13").Value = 1.088414
$ws.Range("J13").Value = 1.565037
$ws.Range("K13").Value = 2.140165
$ws.Range("L13").Value = 2.801308
$ws.Range("M13").Value = 3.570284
$ws.Range("N13").Value = 4.426456
$ws.Range("O13").Value = 5.433976
$ws.Range("P13").Value = 6.631990999999999
$ws.Range("Q13").Value = 7.992939

$ws.Range("C14").Value = 0.013909
$ws.Range("D14").Value = 0.058607
$ws.Range("E14").Value = 0.141561
$ws.Range("F14").Value = 0.276364
$ws.Range("G14").Value = 0.4712370000000001
$ws.Range("H14").Value = 0.7344269999999999
$ws.Range("I14").Value = 1.096746
$ws.Range("J14").Value = 1.569738
$ws.Range("K14").Value = 2.121854
$ws.Range("L14").Value = 2.783945
$ws.Range("M14").Value = 3.555473999999999
$ws.Range("N14").Value = 4.394168000000001
$ws.Range("O14").Value = 5.433046
$ws.Range("P14").Value = 6.586393000000001
$ws.Range("Q14").Value = 7.957806

$ws.Range("C15").Value = 0.014066
$ws.Range("D15").Value = 0.05985500000000001
$ws.Range("E15").Value = 0.145239
$ws.Range("F15").Value = 0.281807
$ws.Range("G15").Value = 0.476622
$ws.Range("H15").Value = 0.7428239999999999
$ws.Range("I15").Value = 1.126362
$ws.Range("J15").Value = 1.601466
$ws.Range("K15").Value = 2.170681
$ws.Range("L15").Value = 2.846819
$ws.Range("M15").Value = 3.631046
$ws.Range("N15").Value = 4.494356000000001
$ws.Range("O15").Value = 5.54616
$ws.Range("P15").Value = 6.726816000000001
$ws.Range("Q15").Value = 8.10455

$ws.Range("C20").Value = 0.029756
$ws.Range("D20").Value = 0.133935
$ws.Range("E20").Value = 0.361437
$ws.Range("F20").Value = 0.7646379999999999
$ws.Range("G20").Value = 1.357092
$ws.Range("H20").Value = 2.187737
$ws.Range("I20").Value = 3.293286
$ws.Range("J20").Value = 4.740823000000001
$ws.Range("K20").Value = 6.539346
$ws.Range("L20").Value = 8.743188999999997
$ws.Range("M20").Value = 11.358781
$ws.Range("N20").Value = 14.370211
$ws.Range("O20").Value = 17.920053
$ws.Range("P20").Value = 22.218637
$ws.Range("Q20").Value = 29.049141

$ws.Range("C21").Value = 0.03153499999999999
$ws.Range("D21").Value = 0.136268
$ws.Range("E21").Value = 0.363836
$ws.Range("F21").Value = 0.772818
$ws.Range("G21").Value = 1.387557
$ws.Range("H21").Value = 2.231125
$ws.Range("I21").Value = 3.309132
$ws.Range("J21").Value = 4.802413
$ws.Range("K21").Value = 6.590045000000001
$ws.Range("L21").Value = 8.767271
$ws.Range("M21").Value = 11.412014
$ws.Range("N21").Value = 14.415702
$ws.Range("O21").Value = 18.048974
$ws.Range("P21").Value = 22.302494
$ws.Range("Q21").Value = 29.116331

$ws.Range("C22").Value = 0.029429
$ws.Range("D22").Value = 0.13759
$ws.Range("E22").Value = 0.364102
$ws.Range("F22").Value = 0.777284
$ws.Range("G22").Value = 1.386954
$ws.Range("H22").Value = 2.220994999999999
$ws.Range("I22").Value = 3.323858
$ws.Range("J22").Value = 4.760835
$ws.Range("K22").Value = 6.57149
$ws.Range("L22").Value = 8.805575999999999
$ws.Range("M22").Value = 11.438485
$ws.Range("N22").Value = 14.377479
$ws.Range("O22").Value = 17.987201
$ws.Range("P22").Value = 22.155565
$ws.Range("Q22").Value = 28.877955

$ws.Range("C23").Value = 0.02991899999999999
$ws.Range("D23").Value = 0.139802
$ws.Range("E23").Value = 0.3741420000000001
$ws.Range("F23").Value = 0.785755
$ws.Range("G23").Value = 1.416995
$ws.Range("H23").Value = 2.269072
$ws.Range("I23").Value = 3.375549
$ws.Range("J23").Value = 4.830209999999999
$ws.Range("K23").Value = 6.690988
$ws.Range("L23").Value = 8.916752999999998
$ws.Range("M23").Value = 11.57588
$ws.Range("N23").Value = 14.547477
$ws.Range("O23").Value = 18.204949
$ws.Range("P23").Value = 22.404405
$ws.Range("Q23").Value = 29.416086

Write-Host "Edit complete"
